$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}

$ws.Cells.Item(5, 3).ClearContents()
$ws.Cells.Item(5, 4).ClearContents()
